$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6499458.5
$ws.Range("J51").Value = 11912262
$ws.Range("L51").Value = 11912262
$ws.Range("N51").Value = -11913230
$ws.Range("H121").Value = 1971.6604
$ws.Range("J121").Value = 1990.3654
$ws.Range("L121").Value = 5971.0962
$ws.Range("N121").Value = -9465.0962
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 9000
$ws.Range("M125").ClearContents()
$ws.Range("N125").Value = -13920
$ws.Range("H131").Value = 11681.363
$ws.Range("I131").Value = 4095
$ws.Range("J131").Value = 12440
$ws.Range("K131").Value = 12285
$ws.Range("L131").Value = 37320
$ws.Range("M131").Value = -7245
$ws.Range("N131").Value = -47400
$ws.Range("H132").Value = 4020.3635
$ws.Range("I132").Value = 3595.2307
$ws.Range("J132").Value = 7336.4
$ws.Range("K132").Value = 10785.6921
$ws.Range("L132").Value = 22009.2
$ws.Range("M132").Value = -8255.6921
$ws.Range("N132").Value = -27069.2
$ws.Range("H137").Value = 3422.625
$ws.Range("I137").Value = 2848.5
$ws.Range("K137").Value = 8545.5
$ws.Range("M137").Value = -5995.5
$ws.Range("H141").Value = 3950
$ws.Range("I141").Value = 1900
$ws.Range("J141").Value = 4633.3335
$ws.Range("K141").Value = 5700
$ws.Range("L141").Value = 13900.0005
$ws.Range("M141").Value = -520
$ws.Range("N141").Value = -24260.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2679.1345
$ws.Range("I32").Value = 2652.5
$ws.Range("J32").Value = 2998.75
$ws.Range("K32").Value = 2652.5
$ws.Range("L32").Value = 2998.75
$ws.Range("M32").Value = -2365.5
$ws.Range("N32").Value = -3572.75
$ws.Range("H45").Value = 1846.6666
$ws.Range("I45").Value = 1444.4
$ws.Range("K45").Value = 1444.4
$ws.Range("M45").Value = -1067.4
$ws.Range("H132").Value = 90911230
$ws.Range("I132").Value = 90911230
$ws.Range("K132").Value = 272733690
$ws.Range("M132").Value = -272731160

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H40").Value = 65999.5
$ws.Range("J40").Value = 65999.5
$ws.Range("L40").Value = 65999.5
$ws.Range("N40").Value = -66529.5
$ws.Range("H105").Value = 1844.35
$ws.Range("J105").Value = 2399.5
$ws.Range("L105").Value = 2399.5
$ws.Range("N105").Value = -5893.5
$ws.Range("H134").Value = 1307.4231
$ws.Range("I134").Value = 1173
$ws.Range("K134").Value = 3519
$ws.Range("M134").Value = -984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3099.4614
$ws.Range("I31").Value = 2165.875
$ws.Range("J31").Value = 4593.2
$ws.Range("K31").Value = 2165.875
$ws.Range("L31").Value = 4593.2
$ws.Range("M31").Value = -1870.875
$ws.Range("N31").Value = -5183.2
$ws.Range("H32").Value = 7999.5
$ws.Range("I32").Value = 7999.5
$ws.Range("K32").Value = 7999.5
$ws.Range("M32").Value = -7683.5
$ws.Range("H34").Value = 3099.4614
$ws.Range("I34").Value = 2165.875
$ws.Range("J34").Value = 4593.2
$ws.Range("K34").Value = 2165.875
$ws.Range("L34").Value = 4593.2
$ws.Range("M34").Value = -1963.875
$ws.Range("N34").Value = -4997.2
$ws.Range("H39").Value = 2961.25
$ws.Range("I39").Value = 2961.25
$ws.Range("K39").Value = 2961.25
$ws.Range("M39").Value = -2570.25
$ws.Range("H49").Value = 2961.25
$ws.Range("I49").Value = 2961.25
$ws.Range("K49").Value = 2961.25
$ws.Range("M49").Value = -2779.25
$ws.Range("H68").Value = 70000
$ws.Range("J68").Value = 70000
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71498
$ws.Range("H71").Value = 70000
$ws.Range("J71").Value = 70000
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -217488
$ws.Range("H132").Value = 1935.8
$ws.Range("I132").Value = 1574.1428
$ws.Range("K132").Value = 4722.428400000001
$ws.Range("M132").Value = -2192.428400000001
$ws.Range("H134").Value = 2552.7273
$ws.Range("I134").Value = 2208.1
$ws.Range("K134").Value = 6624.299999999999
$ws.Range("M134").Value = -4089.299999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 125
$ws.Range("I49").Value = 125
$ws.Range("K49").Value = 375
$ws.Range("M49").Value = -219
$ws.Range("H124").Value = 829.5
$ws.Range("I124").Value = 829.5
$ws.Range("K124").Value = 2488.5
$ws.Range("M124").Value = 2421.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2692.1226
$ws.Range("I132").Value = 2284.1025
$ws.Range("K132").Value = 6852.3075
$ws.Range("M132").Value = -4322.3075

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 33337106
$ws.Range("I7").Value = 62502850
$ws.Range("K7").Value = 62502850
$ws.Range("M7").Value = -62502738
$ws.Range("H126").Value = 33337106
$ws.Range("I126").Value = 62502850
$ws.Range("K126").Value = 187508550
$ws.Range("M126").Value = -187506080
$ws.Range("H132").Value = 5901.8623
$ws.Range("I132").Value = 3688.5625
$ws.Range("J132").Value = 8625.923000000001
$ws.Range("K132").Value = 11065.6875
$ws.Range("L132").Value = 25877.769
$ws.Range("M132").Value = -8535.6875
$ws.Range("N132").Value = -30937.769

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 10008445
$ws.Range("J81").Value = 13342396
$ws.Range("L81").Value = 26684792
$ws.Range("N81").Value = -26686914
$ws.Range("H84").Value = 10008445
$ws.Range("J84").Value = 13342396
$ws.Range("L84").Value = 133423960
$ws.Range("N84").Value = -133434568
$ws.Range("H107").Value = 613
$ws.Range("J107").Value = 360
$ws.Range("L107").Value = 1080
$ws.Range("N107").Value = -4920
$ws.Range("H113").Value = 4334.231
$ws.Range("I113").Value = 149.44444
$ws.Range("J113").Value = 13750
$ws.Range("K113").Value = 448.33332
$ws.Range("L113").Value = 41250
$ws.Range("M113").Value = 1721.66668
$ws.Range("N113").Value = -45590
$ws.Range("H132").Value = 4530.357
$ws.Range("J132").Value = 4468.625
$ws.Range("L132").Value = 13405.875
$ws.Range("N132").Value = -18465.875
$ws.Range("H136").Value = 3637.4211
$ws.Range("I136").Value = 1711.4
$ws.Range("J136").Value = 5777.4443
$ws.Range("K136").Value = 5134.200000000001
$ws.Range("L136").Value = 17332.3329
$ws.Range("M136").Value = -2584.200000000001
$ws.Range("N136").Value = -22432.3329
